$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the EPSG lookup row: reference_system value changes from "EPSG" to "unknown"
$ws.Range("B15").Value = "unknown"

# Move the active cell selection (as reflected in the saved file) near B26
$ws.Range("B26").Select()
